$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the "Conversión del día" note with the new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.0 = 7352.0 pesos`n✅ 7352.0 pesos = 1.99 = 894.17 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate cells N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 500
$ws2.Range("O10").Value = 3676
$ws2.Range("N12").Value = 3699.99
$ws2.Range("O12").Value = 450.001
